$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 467; this shifts the existing rows
# 467..544 down to 468..545 and extends the used range to A1:R545.
$ws.Rows("467:467").Insert()

# Populate the newly inserted row 467 with the new weekly record.
$ws.Cells(467, 1).Value2 = 4
$ws.Cells(467, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells(467, 3).Value2 = "Los Lagos"
$ws.Cells(467, 4).Value2 = 45218
$ws.Cells(467, 5).Value2 = 10
$ws.Cells(467, 6).Value2 = 100112017
$ws.Cells(467, 7).Value2 = "Apio"
$ws.Cells(467, 8).Value2 = "Americana (o)"
$ws.Cells(467, 9).Value2 = "Primera"
$ws.Cells(467, 10).Value2 = 20
$ws.Cells(467, 11).Value2 = 11000
$ws.Cells(467, 12).Value2 = 11000
$ws.Cells(467, 13).Value2 = 11000
$ws.Cells(467, 14).Value2 = "$/docena de matas"
$ws.Cells(467, 15).Value2 = "Región de Coquimbo"
$ws.Cells(467, 16).Value2 = 1833
$ws.Cells(467, 17).Value2 = 6
$ws.Cells(467, 18).Value2 = "Hortaliza"
